$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-sort) values for the 18 data rows (sheet rows 2-19,
# columns A:C) before writing anything back, since the reorder is done in place.
$names = @()
$positions = @()
$teams = @()
for ($r = 2; $r -le 19; $r++) {
    $names += ,$ws.Cells.Item($r, 1).Value()
    $positions += ,$ws.Cells.Item($r, 2).Value()
    $teams += ,$ws.Cells.Item($r, 3).Value()
}

# New order: for each destination data-row index (1-based, 1..18), which
# source data-row index (1-based) supplies its values.
$order = @(14, 13, 3, 6, 7, 8, 11, 12, 18, 5, 4, 2, 1, 10, 15, 16, 17, 9)

for ($i = 0; $i -lt $order.Length; $i++) {
    $src = $order[$i] - 1
    $destRow = $i + 2
    $ws.Cells.Item($destRow, 1).Value = $names[$src]
    $ws.Cells.Item($destRow, 2).Value = $positions[$src]
    $ws.Cells.Item($destRow, 3).Value = $teams[$src]
}
